$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain text in this sheet (e.g. "62.890.53",
# "5.10", "0.0560" -- note the thousands-style dots and the significant
# trailing zeros). If we assign the new values directly, Excel helpfully
# "fixes" anything that looks numeric into a real number (dropping trailing
# zeros, flipping "5.10" into 5.1, etc.), which does not match the original
# text formatting. Switching the whole Price column to the Text number
# format first keeps every assignment literal; resetting the range back to
# the Normal style afterwards keeps the cells on the same (default) style
# they started on.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "62.890.53"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").Value = "2.678.97"
$ws.Range("E3").Value = "  -2.38%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "549.46"
$ws.Range("E5").Value = "  -4.27%  "
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("E9").Value = "  -4.00%  "
$ws.Range("E10").Value = "  -2.95%  "
$ws.Range("E11").Value = "  -4.58%  "
$ws.Range("D12").Value = "5.10"
$ws.Range("E12").Value = "  -11.59%  "
$ws.Range("D13").Value = "3.153.39"
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("D14").Value = "25.95"
$ws.Range("E14").Value = "  -3.65%  "
$ws.Range("D15").Value = "62.761.50"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("D17").Value = "2.681.12"
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("D18").Value = "11.85"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").Value = "4.56"
$ws.Range("E19").Value = "  -5.39%  "
$ws.Range("D20").Value = "342.45"
$ws.Range("E20").Value = "  -3.85%  "
$ws.Range("E21").Value = "  -4.72%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  -4.02%  "
$ws.Range("D24").Value = "63.29"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "8.08"
$ws.Range("E27").Value = "  -4.96%  "
$ws.Range("D28").Value = "0.0₃0853"
$ws.Range("E28").Value = "  -6.57%  "
$ws.Range("D29").Value = "1.93"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("E31").Value = "  -4.21%  "
$ws.Range("D32").Value = "166.71"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("E36").Value = "  -4.84%  "
$ws.Range("E37").Value = "  -4.02%  "
$ws.Range("D38").Value = "336.92"
$ws.Range("E38").Value = "  -3.34%  "
$ws.Range("E39").Value = "  -2.76%  "
$ws.Range("E40").Value = "  -6.98%  "
$ws.Range("E41").Value = "  -2.06%  "
$ws.Range("E42").Value = "  -5.39%  "
$ws.Range("D43").Value = "20.19"
$ws.Range("E43").Value = "  -5.69%  "
$ws.Range("D44").Value = "20.64"
$ws.Range("E44").Value = "  -7.77%  "
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("D46").Value = "0.0560"
$ws.Range("E46").Value = "  -4.85%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").Value = "0.0971"
$ws.Range("E49").Value = "  -3.76%  "
$ws.Range("D50").Value = "128.84"
$ws.Range("E50").Value = "  -4.85%  "
$ws.Range("D51").Value = "2.085.19"
$ws.Range("E51").Value = "  -2.11%  "

$priceCol.Style = "Normal"
